{"js": "// Fill in the \"Summary & Reflection\" closing paragraphs, which previously\n// were five blank paragraphs (carrying only bold paragraph-mark formatting)\n// right before the \"References\" heading, with the final reflection text,\n// and normalize the extra trailing blank paragraphs.\n\nconst body = context.document.body;\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph (\"Each software testing technique...\") that\n// immediately precedes the blank paragraphs - its paragraph-mark formatting\n// (rFonts eastAsiaTheme) is what the new paragraphs should inherit.\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Each software testing technique described above\") === 0) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the anchor paragraph.\");\n}\n\n// The blank paragraphs run from anchorIndex+1 up to (not including) the\n// \"References\" heading paragraph.\nlet refIndex = -1;\nfor (let i = anchorIndex + 1; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"References\") {\n    refIndex = i;\n    break;\n  }\n}\nif (refIndex === -1) {\n  throw new Error(\"Could not find the References heading.\");\n}\n\n// Delete all the blank paragraphs between the anchor and References.\nfor (let i = refIndex - 1; i > anchorIndex; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\nconst paraTexts = [\n  \"My overall mindset in testing this project was to employ a high degree of caution and critical analysis to ensure that requirements were verified and validated. For example, I wrote all requirements down before designing tests and ensured that the tests can catch different failure points. Furthermore, the cumulative testing strategy across multiple files ensured that test coverage was high overall. Taking note of the relationships within and between different files and system components, in my experience, helped in the design of effective and versatile tests.\",\n  \"I attempted to limit bias in my review of the code by writing down requirements and critically analyzing how to design test cases that account for different failure points. For instance, I designed the tests based on the requirements, rather than on what cases would be easiest to write. I attempted to think from the user\\u2019s perspective of potential ways the code could be broken and remedy them. Overall, bias could be a very pertinent concern as a software developer; because bias can lead to errors and defects being missed, ensuring that development and testing are done objectively can increase the chance of product success and a reliable reputation.\",\n  \"Overall, it is important to keep a commitment to quality and not cut corners in the development or testing process because software failure can have consequences at an individual and societal level. For an individual user, software failure can result in program crashes, critical data loss, or inconvenient user experience. Software failure can have massive financial costs, too; problems in the manufacturing of the Airbus A380 aircraft, for example, resulted in direct or indirect costs of approximately 6.1 billion US dollars (Hambling et al., 2019, p. 109). These examples illustrate the importance of keeping a commitment to software quality. Going forward, I plan to avoid technical debt in this field of work by critically analyzing problems, collaborating with my peers, and maintaining a commitment to lifelong learning.\",\n  \"\",\n  \"\",\n  \"\",\n  \"\"\n];\n\n// First-line indent in points: the first paragraph uses 36pt (720 twips /\n// 0.5in); the rest use 18pt (360 twips / 0.25in).\nconst indents = [36, 18, 18, 18, 18, 18, 18];\n\nlet previous = paragraphs.items[anchorIndex];\nfor (let i = 0; i < paraTexts.length; i++) {\n  const newPara = previous.insertParagraph(\"\", Word.InsertLocation.after);\n  await context.sync();\n  if (paraTexts[i] !== \"\") {\n    newPara.insertText(paraTexts[i], Word.InsertLocation.replace);\n  }\n  newPara.firstLineIndent = indents[i];\n  await context.sync();\n  previous = newPara;\n}\n", "ps1": "# Fill in the \"Summary & Reflection\" closing paragraphs, which previously\n# were five blank paragraphs (carrying only bold paragraph-mark formatting)\n# right before the \"References\" heading, with the final reflection text,\n# and normalize the extra trailing blank paragraphs.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"Each software testing technique...\") that\n# immediately precedes the blank paragraphs - its paragraph-mark formatting\n# (rFonts eastAsiaTheme) is what the new paragraphs should inherit.\n$anchorIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"Each software testing technique described above\")) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq 0) {\n    throw \"Could not find the anchor paragraph.\"\n}\n\n# The blank paragraphs run from anchorIndex+1 up to (not including) the\n# \"References\" heading paragraph.\n$refIndex = 0\nfor ($i = $anchorIndex + 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq \"References\") {\n        $refIndex = $i\n        break\n    }\n}\nif ($refIndex -eq 0) {\n    throw \"Could not find the References heading.\"\n}\n\n# Delete all the blank paragraphs between the anchor and References.\nfor ($i = $refIndex - 1; $i -gt $anchorIndex; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n$rsquo = [char]0x2019\n\n$para1 = \"My overall mindset in testing this project was to employ a high degree of caution and critical analysis to ensure that requirements were verified and validated. For example, I wrote all requirements down before designing tests and ensured that the tests can catch different failure points. Furthermore, the cumulative testing strategy across multiple files ensured that test coverage was high overall. Taking note of the relationships within and between different files and system components, in my experience, helped in the design of effective and versatile tests.\"\n$para2 = \"I attempted to limit bias in my review of the code by writing down requirements and critically analyzing how to design test cases that account for different failure points. For instance, I designed the tests based on the requirements, rather than on what cases would be easiest to write. I attempted to think from the user${rsquo}s perspective of potential ways the code could be broken and remedy them. Overall, bias could be a very pertinent concern as a software developer; because bias can lead to errors and defects being missed, ensuring that development and testing are done objectively can increase the chance of product success and a reliable reputation.\"\n$para3 = \"Overall, it is important to keep a commitment to quality and not cut corners in the development or testing process because software failure can have consequences at an individual and societal level. For an individual user, software failure can result in program crashes, critical data loss, or inconvenient user experience. Software failure can have massive financial costs, too; problems in the manufacturing of the Airbus A380 aircraft, for example, resulted in direct or indirect costs of approximately 6.1 billion US dollars (Hambling et al., 2019, p. 109). These examples illustrate the importance of keeping a commitment to software quality. Going forward, I plan to avoid technical debt in this field of work by critically analyzing problems, collaborating with my peers, and maintaining a commitment to lifelong learning.\"\n\n$paraTexts = @($para1, $para2, $para3, \"\", \"\", \"\", \"\")\n# First-line indent in points: the first paragraph uses 36pt (720 twips /\n# 0.5in); the rest use 18pt (360 twips / 0.25in).\n$indents = @(36, 18, 18, 18, 18, 18, 18)\n\n$anchor = $d.Paragraphs.Item($anchorIndex)\n$prev = $anchor\nfor ($i = 0; $i -lt $paraTexts.Length; $i++) {\n    $prev.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($anchorIndex + 1 + $i)\n    if ($paraTexts[$i] -ne \"\") {\n        $newPara.Range.Text = $paraTexts[$i]\n    }\n    $newPara.Format.FirstLineIndent = $indents[$i]\n    $prev = $newPara\n}\n"}
